$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D=44475; I="Primera"; J=250; K=1000; L=1200; M=1100; P=367},
    @{Row=3; D=44161; I="Primera"; J=200; K=600; L=700; M=650; P=217},
    @{Row=4; D=44161; I="Segunda"; J=250; K=500; L=600; M=550; P=183},
    @{Row=5; D=44391; I="Primera"; J=250; K=1800; L=2000; M=1900; P=633},
    @{Row=6; D=44333; I="Primera"; J=200; K=1500; L=1700; M=1600; P=533},
    @{Row=7; D=44523; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=8; D=44523; I="Segunda"; J=200; K=1000; L=1100; M=1050; P=350},
    @{Row=9; D=44481; I="Primera"; J=250; K=950; L=1000; M=975; P=325},
    @{Row=10; D=44467; I="Primera"; J=300; K=1000; L=1200; M=1100; P=367},
    @{Row=11; D=44364; I="Primera"; J=200; K=1700; L=1800; M=1750; P=583},
    @{Row=12; D=44364; I="Segunda"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=13; D=44302; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=14; D=44460; I="Primera"; J=250; K=1400; L=1500; M=1450; P=483},
    @{Row=15; D=44385; I="Primera"; J=200; K=2000; L=2300; M=2150; P=717},
    @{Row=16; D=44174; I="Primera"; J=250; K=500; L=600; M=550; P=183},
    @{Row=17; D=44398; I="Primera"; J=300; K=1700; L=1800; M=1750; P=583},
    @{Row=18; D=44295; I="Primera"; J=200; K=1500; L=1800; M=1650; P=550},
    @{Row=19; D=44249; I="Primera"; J=100; K=1500; L=1600; M=1550; P=517},
    @{Row=20; D=44376; I="Primera"; J=340; K=1400; L=1500; M=1471; P=490},
    @{Row=21; D=44327; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=22; D=44431; I="Primera"; J=250; K=1000; L=1300; M=1150; P=383},
    @{Row=23; D=44306; I="Primera"; J=200; K=2400; L=2500; M=2450; P=817},
    @{Row=24; D=44397; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=25; D=44529; I="Primera"; J=200; K=1000; L=1200; M=1100; P=367},
    @{Row=26; D=44417; I="Primera"; J=250; K=1800; L=2000; M=1900; P=633},
    @{Row=27; D=44417; I="Segunda"; J=200; K=1500; L=1600; M=1550; P=517},
    @{Row=28; D=44432; I="Primera"; J=200; K=1200; L=1300; M=1250; P=417},
    @{Row=29; D=44432; I="Segunda"; J=200; K=950; L=1000; M=975; P=325},
    @{Row=30; D=44428; I="Primera"; J=200; K=1500; L=1800; M=1650; P=550},
    @{Row=31; D=44435; I="Primera"; J=450; K=1000; L=1300; M=1194; P=398},
    @{Row=32; D=44435; I="Segunda"; J=200; K=950; L=1000; M=975; P=325},
    @{Row=33; D=44489; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=34; D=44166; I="Primera"; J=250; K=900; L=1000; M=950; P=317},
    @{Row=35; D=44278; I="Primera"; J=140; K=2000; L=2500; M=2250; P=750},
    @{Row=36; D=44278; I="Segunda"; J=200; K=1500; L=1800; M=1650; P=550},
    @{Row=37; D=44280; I="Primera"; J=200; K=1800; L=2000; M=1900; P=633},
    @{Row=38; D=44280; I="Segunda"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=39; D=44300; I="Primera"; J=250; K=1600; L=1800; M=1700; P=567},
    @{Row=40; D=44270; I="Primera"; J=100; K=1800; L=2000; M=1900; P=633},
    @{Row=41; D=44270; I="Segunda"; J=100; K=1200; L=1500; M=1350; P=450},
    @{Row=42; D=44494; I="Primera"; J=200; K=900; L=1000; M=950; P=317},
    @{Row=43; D=44342; I="Primera"; J=200; K=2000; L=2200; M=2100; P=700},
    @{Row=44; D=44224; I="Primera"; J=200; K=1400; L=1500; M=1450; P=483},
    @{Row=45; D=44224; I="Segunda"; J=160; K=1000; L=1200; M=1100; P=367},
    @{Row=46; D=44447; I="Primera"; J=300; K=1100; L=1200; M=1150; P=383}

)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value = $row.D   # D: Fecha
    $ws.Cells.Item($r, 9).Value = $row.I   # I: Calidad
    $ws.Cells.Item($r, 10).Value = $row.J  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row.K  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row.L  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row.M  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row.P  # P: Precio $/Kg
}
